$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new columns at the front (A:B) -> existing A..G (nome..link)
#    shift right to become C..I.
$ws.Columns("A:B").Insert()

# 2) Insert a new row at row 7 -> existing row 7 (Fonte Carregador Automotivo...)
#    shifts down to row 8, leaving row 7 empty for the new "Controle Universal..." entry.
$ws.Rows("7:7").Insert()

# 3) Header row: fill the two new header cells and match the existing header style
#    (bold, centered, bordered) by copying format from the neighboring header cell.
$ws.Range("A1").Value = "data"
$ws.Range("B1").Value = "loja"
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4) Fill the new "data" / "loja" columns for every data row (2-8).
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = "30/07/2024"
    $ws.Cells.Item($r, 2).Value = "eliteautomotiva"
}

# 5) Update the "link" (column I) values for the rows whose URL tracking params changed.
$ws.Range("I2").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-120a-bob-slim-bivolt-cor-preto/p/MLB22144397?pdp_filters=seller_id:209318924#searchVariation=MLB22144397&position=26&search_layout=grid&type=product&tracking_id=735f2990-9aa7-4837-909c-c7d0ac2697ed"
$ws.Range("I3").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-lite-storm-slim-bivolt-cor-azul/p/MLB24154371?pdp_filters=seller_id:209318924#searchVariation=MLB24154371&position=3&search_layout=grid&type=product&tracking_id=2ec3c3cf-a9d0-41be-9bea-cdd3a115893c"
$ws.Range("I4").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-storm-40a-bivolt-12v-cor-preto/p/MLB22569833?pdp_filters=seller_id:209318924#searchVariation=MLB22569833&position=38&search_layout=grid&type=product&tracking_id=2ec3c3cf-a9d0-41be-9bea-cdd3a115893c"
$ws.Range("I5").Value = "https://produto.mercadolivre.com.br/MLB-3254236266-fonte-automotiva-carregador-jfa-storm-40a-bivolt-12v-digital-_JM#position%3D25%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D17255090-5e75-4a35-8f81-6857e60287d2"
$ws.Range("I6").Value = "https://produto.mercadolivre.com.br/MLB-3709911152-fonte-automotiva-bivolt-jfa-storm-light-200-ampere-carregado-_JM#position%3D23%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D5ba70744-cc61-4d93-8bb9-c5bd5f454d9f"
$ws.Range("I8").Value = "https://produto.mercadolivre.com.br/MLB-3244870575-fonte-carregador-automotivo-jfa-bob-storm-120a-slim-bivolt-_JM#position%3D22%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3Da9318cf0-3742-444b-8d05-f700d91ce4a8"

# 6) Populate the brand-new row 7 ("Controle Universal Longa Distancia Jfa K600 ...").
$ws.Range("C7").Value = "Controle Universal Longa Distância Jfa K600 Alcance 600 Mts"
$ws.Range("D7").Value = "K600"
$ws.Range("E7").Value = 78
$ws.Range("F7").Value = "Acima"
$ws.Range("G7").Value = "NA"
$ws.Range("H7").Value = "premium"
$ws.Range("I7").Value = "https://produto.mercadolivre.com.br/MLB-929012432-controle-universal-longa-distncia-jfa-k600-alcance-600-mts-_JM#position%3D35%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D187cd74c-8a96-4719-8ee1-d447db54679c"
